# Insert a new weekly price record as row 115, shifting all rows from the
# old row 115 onward down by one (old row 115 -> new row 116, ...,
# old row 146 -> new row 147). The worksheet dimension is extended from
# A1:R146 to A1:R147 automatically by the row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 115, pushing existing data down.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new record.
$ws.Range("A115").Value = 7
$ws.Range("B115").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C115").Value = "Ñuble"
$ws.Range("D115").Value = 44511
$ws.Range("E115").Value = 16
$ws.Range("F115").Value = 100112006
$ws.Range("G115").Value = "Repollo"
$ws.Range("H115").Value = "Crespo record"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 300
$ws.Range("K115").Value = 600
$ws.Range("L115").Value = 700
$ws.Range("M115").Value = 650
$ws.Range("N115").Value = '$/unidad'
$ws.Range("O115").Value = "Provincia de Diguillín"
$ws.Range("P115").Value = 650
$ws.Range("Q115").Value = 1
$ws.Range("R115").Value = "Hortaliza"
